# BEAR6-EstimationUX.xlsx — "Fixed bug; worked on ExcelUX"
#
# 1. Rename two worksheets to better reflect their purpose
# 2. Switch the "Estimation start"/"Estimation end" fields over to SMDX-style
#    period strings (quarterly), matching the comments that already explain
#    the expected format
# 3. Populate the previously-empty "Identification horizon" field
# 4. Document the new field with a cell comment
# 5. Update the active cell/selection shown when the sheet is opened

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$wb.Worksheets.Item("Reduced-form estimator").Name = "Reduced-form estimation"
$wb.Worksheets.Item("Structural identifier").Name = "Structural identification"

$wsMeta = $wb.Worksheets.Item("Meta information")

# --- 2. Estimation start / end as SMDX quarterly periods ---------------
$wsMeta.Range("B7").Value = "1975-Q1"
$wsMeta.Range("B8").Value = "2014-Q4"

# --- 3. Identification horizon value -----------------------------------
$wsMeta.Range("B9").Value = 20

# --- 4. Explain the identification horizon field with a comment --------
$comment = $wsMeta.Range("B9").AddComment("Enter an identification horizon as the number of periods for which the impulse responses and FEVDs will be calculated, including in the structural identification step`n")

# --- 5. Move the shown selection to A6 ----------------------------------
$wsMeta.Range("A6").Select() | Out-Null
